# Populate the "Results" sheet with the API call log rows (headers + 3 data rows)
# and make sure the workbook recalculates fully the next time it's opened.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) -------------------------------------------------
$ws.Range("A1").Value = "Operation"
$ws.Range("B1").Value = "ID"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Year"
$ws.Range("E1").Value = "Price"
$ws.Range("F1").Value = "CPU"
$ws.Range("G1").Value = "HardDisk"
$ws.Range("H1").Value = "CreatedAt"
$ws.Range("I1").Value = "Timestamp"

# --- Row 2: POST response -----------------------------------------------
$ws.Range("A2").Value = "POST"
$ws.Range("B2").Value = "ff8081819782e69e0199083b378715dd"
$ws.Range("C2").Value = "Apple MacBook Pro 18 Max"
$ws.Range("D2").Value = 2021
$ws.Range("E2").Value = 189.99
$ws.Range("F2").Value = "Intel Core i8"
$ws.Range("G2").Value = "4 TB"
$ws.Range("H2").Value = "2025-09-02T02:22:02.631+00:00"
$ws.Range("I2").Value = "2025-09-02 07:52:02"

# --- Row 6: GET response for the second (unrelated) record --------------
$ws.Range("A6").Value = "GET"
$ws.Range("B6").Value = "ff8081819782e69e0199083a69e415da"
$ws.Range("C6").Value = "Apple MacBook Pro 18 Max"
$ws.Range("D6").Value = 2021
$ws.Range("E6").Value = 189.99
$ws.Range("F6").Value = "Intel Core i8"
$ws.Range("G6").Value = "4 TB"
$ws.Range("I6").Value = "2025-09-02 07:51:10"

# --- Row 7: GET response for the record created in row 2 -----------------
$ws.Range("A7").Value = "GET"
$ws.Range("B7").Value = "ff8081819782e69e0199083b378715dd"
$ws.Range("C7").Value = "Apple MacBook Pro 18 Max"
$ws.Range("D7").Value = 2021
$ws.Range("E7").Value = 189.99
$ws.Range("F7").Value = "Intel Core i8"
$ws.Range("G7").Value = "4 TB"
$ws.Range("I7").Value = "2025-09-02 07:52:02"

# Touch the page setup so a (blank) <headerFooter/> element is emitted,
# matching the sheet's canonical save shape.
$ws.PageSetup.CenterHeader = ""

# Force a full recalculation the next time the workbook is opened.
$wb.Application.CalculateFull()
